$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.222.84'
$ws.Range("E2").Value = '  -0.53%  '

$ws.Range("D3").Value = '1.909.60'
$ws.Range("E3").Value = '  -1.37%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7332'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.83'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.69%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.35%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3129'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.86'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06914'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7771'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.63%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07990'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.24%  '

$ws.Range("D13").Value = '1.883.35'
$ws.Range("E13").Value = '  -2.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.259'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.40'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.78%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '30.175.25'
$ws.Range("E16").Value = '  -0.70%  '

$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.21'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.850'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.81'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -6.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007792'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.66%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.003'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.33%  '

$ws.Range("D22").Value = '2.143.50'
$ws.Range("E22").Value = '  -2.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.005'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.722'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.406'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.05'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.03'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1268'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.82%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.088'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -8.78%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.548'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.34%  '

$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.346'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.307'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.85%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.083'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05163'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.283'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.43%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7433'
$ws.Range("D36").ClearFormats()

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.758'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01937'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.82%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.797'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.365'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.65'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.60%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4445'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.58%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.932'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.003'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.25%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8352'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.606'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.31%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.06'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.742'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.18%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.53'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '944.63'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.20%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1187'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.53%  '
